# Apply "code refactoring and loan accounting and charges added" changes
# to the Loan Product workbook.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # ProductLoan_Input
$ws2 = $wb.Worksheets.Item(2)   # ProductLoan_Output

# ------------------------------------------------------------------
# 1. Update existing values on the input sheet
# ------------------------------------------------------------------

# Product name: 390-... -> 486-... naming
$ws1.Range("B1").Value = "486-RBI-EPP-DB-DL-REC-NON-RNI-CTPD-DL-MD-TR-1-EarlyRePayment"

# Short name: was the text "kar3", now the numeric code 486
$ws1.Range("B3").Value = 486

# Nominal interest rate default: 12 -> 1
$ws1.Range("B11").Value = 1

# Maximum allowed outstanding balance: 5000 -> 10000
$ws1.Range("B28").Value = 10000

# ------------------------------------------------------------------
# 2. Add the new loan accounting / charges rows (31-42)
# ------------------------------------------------------------------

# Copy the look & feel (fill/font) of an existing label/value row pair
# down onto the new rows before filling in their content.
$ws1.Range("A10:B10").Copy()
$ws1.Range("A31:B42").PasteSpecial(-4122)

$labels = @(
  "fundsource",
  "loanprotfolio",
  "interestreceivable",
  "penaltiesreceivable",
  "transferinsuspense",
  "feesreceivable",
  "incomefrominterest",
  "incomefrompenalties",
  "incomefromfees",
  "incomefromrecoveryrepayments",
  "loseswrittenoff",
  "overpaymentliability"
)

$values = @(
  "Cash",
  "Loan portfolio ",
  "Interest Receivable ",
  "Penalties Receivable ",
  "Transfer in Suspence ",
  "Fees Receivable",
  "Income from interest",
  "Income from penalties",
  "Income from fees",
  "Income from recovery repayments",
  "Losses Writtenoff ",
  "Overpayment Liability"
)

for ($i = 0; $i -lt $labels.Length; $i++) {
  $row = 31 + $i
  $ws1.Cells.Item($row, 1).Value = $labels[$i]
  $ws1.Cells.Item($row, 2).Value = $values[$i]
}

# ------------------------------------------------------------------
# 3. Update the output sheet's product name reference
# ------------------------------------------------------------------

$ws2.Range("B1").Value = "486-RBI-EPP-DB-DL-REC-NON-RNI-CTPD-DL-MD-TR-1-EarlyRePayment"

# ------------------------------------------------------------------
# 4. Restore view/selection state for both sheets
# ------------------------------------------------------------------

$ws2.Activate()
$ws2.Range("B1").Select()

$ws1.Activate()
$ws1.Range("B11").Select()
